$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

# Row 11: Galo Futsal x Vasquinho -> 1x4, Finalizado
$ws.Range("E11").Value = "1x4"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = "Finalizado"

# Row 12: Lagoa Verde x Baixa Grande -> 3x1, Finalizado
$ws.Range("E12").Value = "3x1"
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = "Finalizado"

# Update selection to match diff
$ws.Range("L13").Select()
